# Trade #16 closed at 2026-02-17 20:52:48 - unknown UNKNOWN +0.000%
#
# This script:
#  1. Updates the "Summary" sheet aggregate stats.
#  2. Updates the "Strategy Status" sheet row for MarketMaking.
#  3. Closes trade #44 (row 45 on "All Trades", row 12 on "MarketMaking")
#     as an early_exit loss.
#  4. Appends a brand-new open trade #77 to both "All Trades" (row 78)
#     and "MarketMaking" (row 45).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.35
$summary.Range("B4").Value = 0.14
$summary.Range("B5").Value = 0.06
$summary.Range("B6").Value = 44
$summary.Range("B8").Value = 19
$summary.Range("B9").Value = 43.18

# ---------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.35
$status.Range("D5").Value = 11
$status.Range("E5").Value = 0.03
$status.Range("F5").Value = 0.35
$status.Range("G5").Value = 45.45

# ---------------------------------------------------------------
# 3) All Trades sheet - close trade #44 (row 45)
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(45, 7).Value = 0.18
$allTrades.Cells.Item(45, 8).Value = "CLOSED"
$allTrades.Cells.Item(45, 9).Value = -21.796
$allTrades.Cells.Item(45, 10).Value = -0.05
$allTrades.Cells.Item(45, 11).Value = 100.35
$allTrades.Cells.Item(45, 12).Value = "early_exit"
$allTrades.Cells.Item(45, 13).Value = 0.12

# ---------------------------------------------------------------
# 4) MarketMaking sheet - close trade #44 (row 12)
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(12, 7).Value = 0.18
$marketMaking.Cells.Item(12, 8).Value = "CLOSED"
$marketMaking.Cells.Item(12, 9).Value = -21.796
$marketMaking.Cells.Item(12, 10).Value = -0.05
$marketMaking.Cells.Item(12, 11).Value = 100.35
$marketMaking.Cells.Item(12, 16).Value = "early_exit"
$marketMaking.Cells.Item(12, 17).Value = 0.12

# ---------------------------------------------------------------
# 5) All Trades sheet - append new trade #77 (row 78)
# ---------------------------------------------------------------
$allTrades.Cells.Item(78, 1).Value = 77
$allTrades.Cells.Item(78, 2).NumberFormat = "@"
$allTrades.Cells.Item(78, 2).Value = "2026-02-17"
$allTrades.Cells.Item(78, 3).Value = "20:52:42"
$allTrades.Cells.Item(78, 4).Value = "MarketMaking"
$allTrades.Cells.Item(78, 5).Value = "UP"
$allTrades.Cells.Item(78, 6).Value = 0.230167
$allTrades.Cells.Item(78, 8).Value = "OPEN"
$allTrades.Cells.Item(78, 9).Value = 0
$allTrades.Cells.Item(78, 10).Value = 0
$allTrades.Cells.Item(78, 11).Value = 100.4025618338692
$allTrades.Cells.Item(78, 13).Value = 0
$allTrades.Cells.Item(78, 14).Value = 0
$allTrades.Cells.Item(78, 15).Value = 0
$allTrades.Cells.Item(78, 16).Value = 0.6
$allTrades.Cells.Item(78, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------
# 6) MarketMaking sheet - append new trade #77 (row 45)
# ---------------------------------------------------------------
$marketMaking.Cells.Item(45, 1).Value = 77
$marketMaking.Cells.Item(45, 2).NumberFormat = "@"
$marketMaking.Cells.Item(45, 2).Value = "2026-02-17"
$marketMaking.Cells.Item(45, 3).Value = "20:52:42"
$marketMaking.Cells.Item(45, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(45, 5).Value = "UP"
$marketMaking.Cells.Item(45, 6).Value = 0.230167
$marketMaking.Cells.Item(45, 8).Value = "OPEN"
$marketMaking.Cells.Item(45, 9).Value = 0
$marketMaking.Cells.Item(45, 10).Value = 0
$marketMaking.Cells.Item(45, 11).Value = 100.4025618338692
$marketMaking.Cells.Item(45, 12).Value = 0
$marketMaking.Cells.Item(45, 13).Value = 0
$marketMaking.Cells.Item(45, 14).Value = 0.6
$marketMaking.Cells.Item(45, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(45, 17).Value = 0
